$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 was previously blank (data started at row 2, the header row).
# Write the same header labels into row 1 (without shifting any existing rows).
$headers = @("employee_id","tax_id","firstname","lastname","salary","element1","element2","element3","element4","element5","element6","element7","element8","element9","element10")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Match the bold style used by the existing header row (row 2).
$ws.Range("A1:O1").Font.Bold = $true

# Update the active selection to Q3, as in the edited workbook.
$ws.Range("Q3").Select()
